$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:G to B:H.
$ws.Columns.Item(1).Insert()

# Number of data rows (rows 2..21 in the new layout correspond to the
# original data rows).
$lastRow = 21

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Give the new index column the same header-row style (bold, centered,
# bordered) that the other header cells already use, without touching
# the numeric values we just wrote.
$ws.Range("B1").Copy()
$ws.Range("A2:A$lastRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
